# Remove the old "SaveCommand" actor lane (Rectangle 65 "Logic" lane,
# its lifeline/activation, the execute() call, :SaveCommand box, its
# lifeline/return arrows, the X-destroy marker, the :CommandResult
# return group, and the saveCommand() label) from the sequence diagram
# on the single slide, and nudge the remaining "Logic -> :SaveSequence"
# return arrow (Id 50) to its new anchor point.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shapes.Item() indices (1-based) of every shape slated for removal,
# paired with the shape Id we expect to find there (taken from the
# slide in its original order). Deleting from the highest index down
# keeps the remaining indices stable while we work.
$toRemove = @(
    @{ Index = 31; Id = 57 },   # TextBox 56 ("saveCommand()")
    @{ Index = 26; Id = 232 },  # Straight Arrow Connector 231
    @{ Index = 25; Id = 229 },  # Straight Arrow Connector 228
    @{ Index = 24; Id = 223 },  # Group 222 (:CommandResult group)
    @{ Index = 20; Id = 94 },   # TextBox 93 ("X")
    @{ Index = 19; Id = 65 },   # Rectangle 64
    @{ Index = 16; Id = 69 },   # Straight Arrow Connector 68
    @{ Index = 15; Id = 66 },   # Straight Arrow Connector 65
    @{ Index = 14; Id = 19 },   # Rectangle 62 (":SaveCommand")
    @{ Index = 8;  Id = 34 },   # Straight Arrow Connector 33
    @{ Index = 7;  Id = 29 },   # TextBox 28 ("execute()")
    @{ Index = 6;  Id = 28 },   # Straight Arrow Connector 27
    @{ Index = 5;  Id = 21 },   # Rectangle 20
    @{ Index = 4;  Id = 20 },   # Straight Connector 19
    @{ Index = 3;  Id = 81 }    # Rectangle 65 ("Logic")
)

foreach ($entry in $toRemove) {
    $sh = $s.Shapes.Item($entry.Index)
    if ($sh.Id -eq $entry.Id) {
        $sh.Delete()
    } else {
        # Fall back to an Id search if the slide wasn't in the
        # expected original order for some reason.
        $n = $s.Shapes.Count
        for ($i = 1; $i -le $n; $i++) {
            $candidate = $s.Shapes.Item($i)
            if ($candidate.Id -eq $entry.Id) {
                $candidate.Delete()
                break
            }
        }
    }
}

# Re-anchor "Straight Arrow Connector 49" (Id 50) now that the lane
# above it is gone; look it up by its stable Id since the collection
# was just re-indexed by the deletions above.
$n = $s.Shapes.Count
for ($i = 1; $i -le $n; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 50) {
        $sh.Left = 347.82511905511814
        $sh.Top = 252.851968503937
        $sh.Width = 231.3623622047244
        $sh.Height = 0.0
        break
    }
}
